# MidAI.xlsx edit: add "Romania" (edge list) and "Heuristic" (straight-line
# distance to Bucharest) sheets after the existing "Problem1" sheet, and
# tweak a few view/format properties on all three sheets.

$wb = $excel.ActiveWorkbook
$sheets = $wb.Worksheets

# ---------------------------------------------------------------------
# 1. Create the two new worksheets, in tab order Problem1, Romania,
#    Heuristic.
# ---------------------------------------------------------------------
$problem1 = $sheets.Item(1)

$romania = $sheets.Add($null, $problem1)
$romania.Name = "Romania"

$heuristic = $sheets.Add($null, $romania)
$heuristic.Name = "Heuristic"

# ---------------------------------------------------------------------
# 2. Romania sheet: VerticeA / VerticeB / Weight edge list for the
#    classic Romania road-map graph (AIMA).
# ---------------------------------------------------------------------
$romania.Cells.Item(1, 1).Value = "VerticeA"
$romania.Cells.Item(1, 2).Value = "VerticeB"
$romania.Cells.Item(1, 3).Value = "Weight"

$edges = @(
  @("Oradea","Zerind",71),
  @("Oradea","Sibiu",151),
  @("Zerind","Arad",75),
  @("Arad","Timisoara",118),
  @("Arad","Sibiu",140),
  @("Lugoj","Timisoara",111),
  @("Lugoj","Mehadia",70),
  @("Dobreta","Mehadia",75),
  @("Dobreta","Craiova",120),
  @("Fagaras","Sibiu",99),
  @("Rimnicu Vilcea","Sibiu",80),
  @("Pitesti","Craiova",138),
  @("Rimnicu Vilcea","Pitesti",97),
  @("Rimnicu Vilcea","Craiova",146),
  @("Fagaras","Bucharest",211),
  @("Pitesti","Bucharest",101),
  @("Giurgiu","Bucharest",90),
  @("Urziceni","Bucharest",85),
  @("Urziceni","Hirsova",98),
  @("Eforie","Hirsova",86),
  @("Urziceni","Vaslui",142),
  @("Iasi","Vaslui",92),
  @("Iasi","Neamt",87)
)

$r = 2
foreach ($edge in $edges) {
  if ($r -ne 18) {
    $romania.Cells.Item($r, 1).Value = $edge[0]
  }
  $romania.Cells.Item($r, 2).Value = $edge[1]
  $romania.Cells.Item($r, 3).Value = $edge[2]
  $r = $r + 1
}

# ---------------------------------------------------------------------
# 3. Heuristic sheet: straight-line distance to Bucharest.
# ---------------------------------------------------------------------
$heuristic.Cells.Item(1, 1).Value = "To Bucharest"
$heuristic.Cells.Item(1, 2).Value = "Straight-line distance"

$heur = @(
  @("Arad",366),
  @("Bucharest",0),
  @("Cralova",160),
  @("Dobreta",242),
  @("Eforie",161),
  @("Fagaras",178),
  @("Giurgiu",77),
  @("Hirsova",151),
  @("Iasi",226),
  @("Lugoj",244),
  @("Mehadia",241),
  @("Neamt",234),
  @("Oradea",380),
  @("Pitesti",98),
  @("Rimnicu Vilcea",193),
  @("Sibiu",253),
  @("Timisoara",329),
  @("Urziceni",80),
  @("Vaslui",199),
  @("Zerind",374)
)

$r = 2
foreach ($h in $heur) {
  $heuristic.Cells.Item($r, 1).Value = $h[0]
  $heuristic.Cells.Item($r, 2).Value = $h[1]
  $r = $r + 1
}

# ---------------------------------------------------------------------
# 4. Go back and fill the one cell that was skipped (Romania!A18) last,
#    so it is appended to the shared-string table after the Heuristic
#    sheet's strings -- matching how the source workbook was authored.
# ---------------------------------------------------------------------
$romania.Cells.Item(18, 1).Value = "Giurgiu"

# ---------------------------------------------------------------------
# 5. Column widths (approximate Excel's "best fit" autosize).
# ---------------------------------------------------------------------
$romania.Columns.Item(1).ColumnWidth = 13.43
$romania.Columns.Item(2).ColumnWidth = 9.0

$heuristic.Columns.Item(1).ColumnWidth = 13.43
$heuristic.Columns.Item(2).ColumnWidth = 19.43

# ---------------------------------------------------------------------
# 6. Selections / view state per sheet.
# ---------------------------------------------------------------------
$problem1.Range("G15:G16").Select()
$romania.Range("A4").Select()
$heuristic.Range("E19").Select()

# ---------------------------------------------------------------------
# 7. Heuristic ends up as the active (selected) tab.
# ---------------------------------------------------------------------
$heuristic.Activate()
